$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.485.56'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3705'
$ws.Range('E7').Value = '  -2.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.89'
$ws.Range('E8').Value = '  +1.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3387'
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07547'
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.142'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.020'
$ws.Range('E14').Value = '  +0.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.957'
$ws.Range('D16').Value = '1.573.50'
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001119'
$ws.Range('E17').Value = '  -1.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.66'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06756'
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.306'
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.23'
$ws.Range('D24').Value = '22.480.30'
$ws.Range('E24').Value = '  +0.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.375'
$ws.Range('E25').Value = '  -1.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.605'
$ws.Range('E26').Value = '  -3.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.05'
$ws.Range('E27').Value = '  -0.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '149.18'
$ws.Range('E28').Value = '  +1.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.055'
$ws.Range('E29').Value = '  +0.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.11'
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('D31').Value = '1.748.51'
$ws.Range('E31').Value = '  +0.49%  '
$ws.Range('E32').Value = '  +7.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.242'
$ws.Range('E33').Value = '  +2.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.014'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.771'
$ws.Range('E35').Value = '  -3.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08353'
$ws.Range('E36').Value = '  -1.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02490'
$ws.Range('E37').Value = '  -0.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.367'
$ws.Range('E38').Value = '  -4.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2303'
$ws.Range('E39').Value = '  +0.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06546'
$ws.Range('E40').Value = '  +0.94%  '
$ws.Range('E41').Value = '  +0.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.35'
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6231'
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.03'
$ws.Range('E45').Value = '  +0.76%  '
$ws.Range('E46').Value = '  +0.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5859'
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.18'
$ws.Range('E48').Value = '  +3.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.074'
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('E50').Value = '  -2.64%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07329'
$ws.Range('E51').Value = '  +0.09%  '
